$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a serial date value (45190 -> 2023-09-21).
# Update every data row (2 through 118) to the new serial date 45192 (2023-09-23).
for ($row = 2; $row -le 118; $row++) {
    $ws.Cells.Item($row, 3).Value = 45192
}
